# Insert a new row at position 89 (shifts existing rows 89-141 down to 90-142)
# and populate it with the new data point, per the authors' commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(89).Insert()

$row = 89
$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 45236
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = 100112026
$ws.Cells.Item($row, 7).Value  = "Haba"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 700
$ws.Cells.Item($row, 11).Value = 8000
$ws.Cells.Item($row, 12).Value = 9000
$ws.Cells.Item($row, 13).Value = 8429
$ws.Cells.Item($row, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región del Maule"
$ws.Cells.Item($row, 16).Value = 337
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
